$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("D2").Value = 0.5489325840217185
$ws.Range("E2").Value = 0.2484647400409127

$ws.Range("D3").Value = 0.5462247212009397
$ws.Range("E3").Value = 0.2420163501452704

$ws.Range("D4").Value = 0.5444788569396176
$ws.Range("E4").Value = 0.2390291355613186
